$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "66.342.66"
$ws.Range("D3").Value = "3.322.23"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5"); $c.NumberFormat = "@"; $c.Value = "588.00"; $c.Style = "Normal"
$ws.Range("E5").Value = "  +2.61%  "
$c = $ws.Range("D6"); $c.NumberFormat = "@"; $c.Value = "183.90"; $c.Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  +8.11%  "
$ws.Range("E9").Value = "  -3.13%  "
$ws.Range("E10").Value = "  +1.99%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "3.901.50"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("D14").Value = "66.351.24"
$ws.Range("E14").Value = "  -0.82%  "
$c = $ws.Range("D15"); $c.NumberFormat = "@"; $c.Value = "26.23"; $c.Style = "Normal"
$ws.Range("E16").Value = "  -2.74%  "
$ws.Range("D17").Value = "3.282.81"
$ws.Range("E17").Value = "  -1.90%  "
$c = $ws.Range("D18"); $c.NumberFormat = "@"; $c.Value = "426.71"; $c.Style = "Normal"
$ws.Range("E18").Value = "  -2.29%  "
$c = $ws.Range("D19"); $c.NumberFormat = "@"; $c.Value = "5.53"; $c.Style = "Normal"
$ws.Range("E19").Value = "  -2.78%  "
$c = $ws.Range("D20"); $c.NumberFormat = "@"; $c.Value = "13.20"; $c.Style = "Normal"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("E21").Value = "  -2.96%  "
$c = $ws.Range("D22"); $c.NumberFormat = "@"; $c.Value = "71.92"; $c.Style = "Normal"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "3.463.66"
$ws.Range("E25").Value = "  -0.79%  "
$c = $ws.Range("D26"); $c.NumberFormat = "@"; $c.Value = "0.515"; $c.Style = "Normal"
$c = $ws.Range("D27"); $c.NumberFormat = "@"; $c.Value = "0.204"; $c.Style = "Normal"
$ws.Range("E27").Value = "  +7.04%  "
$ws.Range("E28").Value = "  -3.94%  "
$ws.Range("E29").Value = "  -1.06%  "
$c = $ws.Range("D30"); $c.NumberFormat = "@"; $c.Value = "0.999"; $c.Style = "Normal"
$ws.Range("E30").Value = "  +0.02%  "
$c = $ws.Range("D31"); $c.NumberFormat = "@"; $c.Value = "1.94"; $c.Style = "Normal"
$ws.Range("E31").Value = "  -1.15%  "
$c = $ws.Range("D32"); $c.NumberFormat = "@"; $c.Value = "22.33"; $c.Style = "Normal"
$ws.Range("E32").Value = "  -2.34%  "
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("E34").Value = "  -2.53%  "
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("E36").Value = "  -3.89%  "
$c = $ws.Range("D37"); $c.NumberFormat = "@"; $c.Value = "159.83"; $c.Style = "Normal"
$ws.Range("E37").Value = "  -1.07%  "
$ws.Range("E38").Value = "  -3.92%  "
$ws.Range("D39").Value = "2.889.70"
$ws.Range("E39").Value = "  +1.98%  "
$c = $ws.Range("D40"); $c.NumberFormat = "@"; $c.Value = "1.80"; $c.Style = "Normal"
$ws.Range("E40").Value = "  -2.53%  "
$c = $ws.Range("D41"); $c.NumberFormat = "@"; $c.Value = "26.47"; $c.Style = "Normal"
$ws.Range("E41").Value = "  -5.06%  "
$c = $ws.Range("D42"); $c.NumberFormat = "@"; $c.Value = "0.765"; $c.Style = "Normal"
$ws.Range("E42").Value = "  -3.38%  "
$c = $ws.Range("D43"); $c.NumberFormat = "@"; $c.Value = "4.32"; $c.Style = "Normal"
$ws.Range("E43").Value = "  -3.00%  "
$ws.Range("E44").Value = "  +0.03%  "
$c = $ws.Range("D45"); $c.NumberFormat = "@"; $c.Value = "0.0665"; $c.Style = "Normal"
$ws.Range("E45").Value = "  -1.25%  "
$c = $ws.Range("D46"); $c.NumberFormat = "@"; $c.Value = "5.93"; $c.Style = "Normal"
$ws.Range("E46").Value = "  -6.57%  "
$ws.Range("E47").Value = "  -2.63%  "
$ws.Range("E48").Value = "  -5.50%  "
$c = $ws.Range("D49"); $c.NumberFormat = "@"; $c.Value = "314.68"; $c.Style = "Normal"
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("E51").Value = "  +5.02%  "
